$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 675
$ws1.Range("F14").Value = 2172
$ws1.Range("F18").Value = 274

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 23
$ws2.Range("F10").Value = 42

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 675
$ws4.Range("F18").Value = 2173
$ws4.Range("F23").Value = 274
$ws4.Range("F25").Value = 23
$ws4.Range("F32").Value = 42
